# Update "想去人数" (F column) counts on several rows across the
# 展览, 演出 and 全部类型 sheets, as reflected in the regenerated
# gh-pages data output.

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (Exhibitions) --
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 202
$ws1.Range("F6").Value  = 1268
$ws1.Range("F14").Value = 278
$ws1.Range("F17").Value = 575
$ws1.Range("F19").Value = 368
$ws1.Range("F21").Value = 875
$ws1.Range("F25").Value = 2707
$ws1.Range("F30").Value = 825
$ws1.Range("F31").Value = 1406
$ws1.Range("F37").Value = 682
$ws1.Range("F38").Value = 708
$ws1.Range("F39").Value = 906

# -- Sheet "演出" (Performances) --
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 34
$ws2.Range("F10").Value = 5
$ws2.Range("F15").Value = 700

# -- Sheet "全部类型" (All types) --
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 202
$ws4.Range("F9").Value  = 1268
$ws4.Range("F16").Value = 278
$ws4.Range("F20").Value = 575
$ws4.Range("F22").Value = 368
$ws4.Range("F26").Value = 2707
$ws4.Range("F35").Value = 825
$ws4.Range("F36").Value = 1406
$ws4.Range("F42").Value = 682
$ws4.Range("F43").Value = 708
$ws4.Range("F44").Value = 906
